# Generate Report for Handback
# Updates the "Ready for handoff" row (row 7) on the zh-cn and de-de
# worksheets now that a (stale) handback has come in for the
# 6c154245-f21e-4e71-a048-4cfbacb74052 source file: fill in the target
# file link, the handback file name, the handback datetime, and the
# "not the latest" error detail.

$wb = $excel.ActiveWorkbook

$targetFileName = "6c154245-f21e-4e71-a048-4cfbacb74052.md"
$targetFileUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cdc62e21d323fc12dc0c282dde7a8af672ceaee/e2e/6c154245-f21e-4e71-a048-4cfbacb74052.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f853a578443f169ca66f9454efcb14e5201eab60/e2e/6c154245-f21e-4e71-a048-4cfbacb74052.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2cdc62e21d323fc12dc0c282dde7a8af672ceaee/e2e/6c154245-f21e-4e71-a048-4cfbacb74052.md."

# Per-language: handback xliff file name + handback datetime.
$zhcn = @{
    Sheet = "zh-cn"
    HandbackFile = "6c154245-f21e-4e71-a048-4cfbacb74052.46d39db84883616797e41a1e74f53d67036d0d03.zh-cn.xlf"
    HandbackDate = "2016-08-28 07:00:59"
}
$dede = @{
    Sheet = "de-de"
    HandbackFile = "6c154245-f21e-4e71-a048-4cfbacb74052.46d39db84883616797e41a1e74f53d67036d0d03.de-de.xlf"
    HandbackDate = "2016-08-28 07:01:11"
}

foreach ($row in @($zhcn, $dede)) {
    $ws = $wb.Worksheets.Item($row.Sheet)

    # Latest Target File (I7) - hyperlinked, same as the source file link.
    $ws.Hyperlinks.Add($ws.Range("I7"), $targetFileUrl, $null, $null, $targetFileName)

    # Latest Handback File (J7)
    $ws.Range("J7").Value = $row.HandbackFile

    # Latest Handback DateTime (K7)
    $ws.Range("K7").Value = $row.HandbackDate

    # Error Detail (P7)
    $ws.Range("P7").Value = $errorDetail
}
